$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 94

# Columns A, B, D hold values that look numeric/date-like ("-548", "8/6/2025",
# "13") but must stay as plain text, matching every other row in this sheet.
# Forcing NumberFormat to Text ("@") before assigning the value stops Excel
# from auto-converting them into a number/date; ClearFormats() afterwards
# drops the now-unneeded explicit cell style so the new row keeps the same
# (default/no-style) look as the rest of the data rows.
$textForceCols = 1, 2, 4

foreach ($col in $textForceCols) {
    $ws.Cells.Item($row, $col).NumberFormat = "@"
}

$ws.Cells.Item($row, 1).Value = "-548"
$ws.Cells.Item($row, 2).Value = "8/6/2025"
$ws.Cells.Item($row, 3).Value = "Sucre 1533"
$ws.Cells.Item($row, 4).Value = "13"
$ws.Cells.Item($row, 5).Value = "Pendiente ADM"
$ws.Cells.Item($row, 6).Value = "Optical Power"
$ws.Cells.Item($row, 7).Value = "Pendiente"
$ws.Cells.Item($row, 8).Value = "Traspasar red a columna de TLC y Retirar columna quebrada"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = "Desmonte"
$ws.Cells.Item($row, 11).Value = "Sin equipos"
$ws.Cells.Item($row, 12).Value = "Pasante"
$ws.Cells.Item($row, 13).Value = -58.44649
$ws.Cells.Item($row, 14).Value = -34.558808
$ws.Cells.Item($row, 15).Value = "Saavedra"
$ws.Cells.Item($row, 16).Value = "Capital Norte"

foreach ($col in $textForceCols) {
    $ws.Cells.Item($row, $col).ClearFormats()
}
